$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value = "tower_group"
for ($r = 2; $r -le 29; $r++) {
    $ws.Cells.Item($r, 7).Value = 1
}

$ws.Range("G12").Select()
